# Update rails to accept chart date range:
# Add two new columns ("start_date" and "end_date", both type "string")
# to the "charts" table schema, directly below "chart_type" and above
# the blank/divider rows that close out that table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "charts" table currently looks like (rows 87-93):
#   87  charts                (header)
#   88  COLUMN | DATA TYPE | VALIDATION
#   89  id          | primary key
#   90  company_id  | foreign key
#   91  chart_type  | string  | ['candle', 'xy', etc]
#   92  (blank divider row)
#   93  (thick-bottom divider row)
#
# We insert two new blank rows right above the current "chart_type" row
# (row 91), which pushes chart_type (and everything below it) down by
# two rows, then fill those two new rows with the new columns.

$ws.Range("B91:D92").Insert(-4121) | Out-Null   # -4121 = xlShiftDown

$ws.Cells.Item(91, 2).Value = "start_date"
$ws.Cells.Item(91, 3).Value = "string"

$ws.Cells.Item(92, 2).Value = "end_date"
$ws.Cells.Item(92, 3).Value = "string"

# Match the formatting of a normal (non-validation) row in this table,
# e.g. row 90 (company_id), which has the plain left/right border
# treatment used throughout these schema tables.
$ws.Range("B90:D90").Copy() | Out-Null
$ws.Range("B91:D91").PasteSpecial(-4122) | Out-Null   # -4122 = xlPasteFormats
$ws.Range("B92:D92").PasteSpecial(-4122) | Out-Null   # -4122 = xlPasteFormats
$excel.CutCopyMode = $false

# Reflect the author's new cursor position after the edit.
$ws.Range("C93").Select() | Out-Null
